$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: length of roads, formatted like D13, value added first so it
# lands at shared-string index 20 (matches target ordering).
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "Длина дорог - roadslen (км) (id8006005)"

# New column E header (row10) "Сельское хозяйство", formatted like D10.
$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "Сельское хозяйство"

# Extend the formatting of column E down through rows 11-13 (no values).
$ws.Range("D11").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$ws.Range("D13").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to match target state
$ws.Range("E17").Select()
